# Recolor the "Minimum Window Substring" / "Valid Anagram" block to green (00B050)
# and, along the way, split the "window" word of the Minimum-Window-Substring URL
# into separate single-letter runs ("w" / "i") the way Word does when text gets
# touched up interactively.

$d = $word.ActiveDocument
$green = 5287936  # 0x00B050 as wdColor (R + G*256 + B*65536)

# --- 1. Line break right before "- Minimum Window Substring - " : 222222 -> 00B050
$rng = $d.Content
$rng.Find.Execute("- Minimum Window Substring - ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$brBeforeMin = $d.Range($rng.Start - 1, $rng.Start)
$brBeforeMin.Font.Color = $green

# --- 2. Split the Minimum Window Substring URL into 4 runs and recolor them green
$urlRng = $d.Content
$urlRng.Find.Execute("https://leetcode.com/problems/minimum-window-substring/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$urlStart = $urlRng.Start

$wChar = $d.Range($urlStart + 38, $urlStart + 39)
$wChar.Bold = 1
$wChar.Bold = 0
$wChar.Font.Color = $green

$iChar = $d.Range($urlStart + 39, $urlStart + 40)
$iChar.Bold = 1
$iChar.Bold = 0
$iChar.Font.Color = $green

# --- 3. Line break right after the hyperlink (before "- Valid Anagram - ") : 222222 -> 00B050
$vaRng = $d.Content
$vaRng.Find.Execute("- Valid Anagram - ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$brBeforeVa = $d.Range($vaRng.Start - 1, $vaRng.Start)
$brBeforeVa.Font.Color = $green

# --- 4. "- Valid Anagram - " text itself : 222222 -> 00B050
$vaRng.Font.Color = $green

# --- 5. Valid Anagram hyperlink URL : 37ACC9 -> 00B050
$vaUrlRng = $d.Content
$vaUrlRng.Find.Execute("https://leetcode.com/problems/valid-anagram/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$vaUrlRng.Font.Color = $green
